$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look numeric, so Excel
# keeps them as literal text (matching the original inline-string cells)
# instead of auto-converting to a number.
foreach ($cellref in @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D16", "D18", "D19", "D20", "D22", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D47", "D48", "D49", "D50", "D51")) {
    $ws.Range($cellref).NumberFormat = "@"
}

$ws.Range("D2").Value = '25.880.81'
$ws.Range("E2").Value = '  -0.20%  '

$ws.Range("D3").Value = '1.631.37'
$ws.Range("E3").Value = '  -0.62%  '

$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.45%  '

$ws.Range("D5").Value = '214.24'
$ws.Range("E5").Value = '  -0.48%  '

$ws.Range("D6").Value = '0.5131'
$ws.Range("E6").Value = '  +1.84%  '

$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  -0.28%  '

$ws.Range("D8").Value = '0.2553'
$ws.Range("E8").Value = '  -0.68%  '

$ws.Range("D9").Value = '0.06338'
$ws.Range("E9").Value = '  -1.21%  '

$ws.Range("D10").Value = '19.42'
$ws.Range("E10").Value = '  -0.91%  '

$ws.Range("D11").Value = '0.07741'
$ws.Range("E11").Value = '  -0.54%  '

$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = '4.258'
$ws.Range("E12").Value = '  -0.26%  '

$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.639.55'
$ws.Range("E13").Value = '  -0.22%  '

$ws.Range("D14").Value = '0.5411'
$ws.Range("E14").Value = '  -0.34%  '

$ws.Range("D15").Value = '0.0₅7718'
$ws.Range("E15").Value = '  -2.66%  '

$ws.Range("D16").Value = '64.03'
$ws.Range("E16").Value = '  -0.69%  '

$ws.Range("D17").Value = '25.884.08'
$ws.Range("E17").Value = '  -0.33%  '

$ws.Range("D18").Value = '1.001'
$ws.Range("E18").Value = '  -0.35%  '

$ws.Range("D19").Value = '195.15'
$ws.Range("E19").Value = '  -2.03%  '

$ws.Range("D20").Value = '4.403'
$ws.Range("E20").Value = '  +0.38%  '

$ws.Range("E21").Value = '  -0.18%  '

$ws.Range("D22").Value = '6.018'
$ws.Range("E22").Value = '  +0.85%  '

$ws.Range("E23").Value = '  -0.05%  '

$ws.Range("D24").Value = '1.853'
$ws.Range("E24").Value = '  -0.92%  '

$ws.Range("D25").Value = '141.20'
$ws.Range("E25").Value = '  +0.04%  '

$ws.Range("D26").Value = '0.1198'
$ws.Range("E26").Value = '  +5.62%  '

$ws.Range("D27").Value = '6.800'
$ws.Range("E27").Value = '  -0.18%  '

$ws.Range("E28").Value = '  -0.95%  '

$ws.Range("D29").Value = '1.232'
$ws.Range("E29").Value = '  -0.78%  '

$ws.Range("D30").Value = '0.04902'
$ws.Range("E30").Value = '  -0.52%  '

$ws.Range("D31").Value = '3.227'
$ws.Range("E31").Value = '  -1.17%  '

$ws.Range("D32").Value = '3.146'
$ws.Range("E32").Value = '  -1.82%  '

$ws.Range("D33").Value = '1.525'
$ws.Range("E33").Value = '  -0.97%  '

$ws.Range("D34").Value = '2.364'
$ws.Range("E34").Value = '  -0.24%  '

$ws.Range("D35").Value = '0.8854'
$ws.Range("E35").Value = '  -0.84%  '

$ws.Range("D36").Value = '2.568'
$ws.Range("E36").Value = '  -1.50%  '

$ws.Range("D37").Value = '1.131.68'
$ws.Range("E37").Value = '  -1.53%  '

$ws.Range("D38").Value = '0.5380'
$ws.Range("E38").Value = '  -3.14%  '

$ws.Range("D39").Value = '0.01547'
$ws.Range("E39").Value = '  -1.47%  '

$ws.Range("D40").Value = '1.002'
$ws.Range("E40").Value = '  -0.28%  '

$ws.Range("D41").Value = '2.533'
$ws.Range("E41").Value = '  -1.05%  '

$ws.Range("D42").Value = '0.8103'
$ws.Range("E42").Value = '  -0.14%  '

$ws.Range("D43").Value = '5.455'
$ws.Range("E43").Value = '  -4.64%  '

$ws.Range("D44").Value = '98.91'
$ws.Range("E44").Value = '  -0.97%  '

$ws.Range("E45").Value = '  +3.15%  '

$ws.Range("D46").Value = '1.769.18'
$ws.Range("E46").Value = '  -0.57%  '

$ws.Range("D47").Value = '0.4519'
$ws.Range("E47").Value = '  -0.06%  '

$ws.Range("D48").Value = '1.004'
$ws.Range("E48").Value = '  -0.04%  '

$ws.Range("D49").Value = '54.47'
$ws.Range("E49").Value = '  -0.45%  '

$ws.Range("D50").Value = '0.05045'
$ws.Range("E50").Value = '  +0.07%  '

$ws.Range("D51").Value = '1.001'
$ws.Range("E51").Value = '  -0.37%  '
